$d = $word.ActiveDocument

# --- 1. Remove the old "_GoBack" bookmark (was a collapsed point right
#        after "...P/E ratio, Beta") so the id=0 slot is free again. ---
$oldBm = $d.Bookmarks.Item("_GoBack")
$oldBm.Delete()

# --- 2. Update the "Trading platform ..." bullet text ---
# "Trading platform for Interest Rate Swaps (debt instruments with $1,000,000 minimum value)"
#   -> "Swap Execution Facility for Interest Rate Swaps (debt instruments with $1,000,000 min. value)"
$d.Content.Find.Execute(
    "Trading platform", $true, $false, $false, $false, $false,
    $true, 1, $false, "Swap Execution Facility", 2) | Out-Null

$d.Content.Find.Execute(
    "minimum value", $true, $false, $false, $false, $false,
    $true, 1, $false, "min. value", 2) | Out-Null

# --- 3. Re-add the "_GoBack" bookmark so it now spans the whole
#        "Swap Execution Facility ..." paragraph (the bullet that was
#        just edited), matching where the author last worked. ---
$p4 = $d.Paragraphs.Item(4)
$p5 = $d.Paragraphs.Item(5)
$bmRange = $d.Range($p4.Range.Start, $p5.Range.Start)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Output "done"
